# adding AE2511 to discrete file
# The workbook has two sheets: "mergedSections" and "fullGrid".
# On "fullGrid", update two filenames in column C:
#   C14: "Create_bioscope_files_2024_Krista.m" -> "Create_biosscope_files_2026.m"
#   C20: "Join_discreteData.R" -> "Join_discreteData_v3.R"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fullGrid")

$ws.Range("C14").Value = "Create_biosscope_files_2026.m"
$ws.Range("C20").Value = "Join_discreteData_v3.R"

$ws.Range("C20").Select()
